# Update cryptocurrency price/volume data in cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, otherwise Excel auto-converts them to
# numeric values (losing trailing zeros / exact formatting).
$forceTextAddrs = @(
    "D5",
    "D6",
    "D11",
    "D12",
    "D16",
    "D17",
    "D18",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D28",
    "D31",
    "D32",
    "D34",
    "D40",
    "D41",
    "D42",
    "D45",
    "D46",
    "D49",
    "D51"
)
foreach ($addr in $forceTextAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Set the new values for the force-text cells
$ws.Range("D5").Value = "553.17"
$ws.Range("D6").Value = "159.42"
$ws.Range("D11").Value = "0.328"
$ws.Range("D12").Value = "4.79"
$ws.Range("D16").Value = "10.33"
$ws.Range("D17").Value = "333.65"
$ws.Range("D18").Value = "6.82"
$ws.Range("D22").Value = "66.13"
$ws.Range("D23").Value = "3.62"
$ws.Range("D24").Value = "8.09"
$ws.Range("D26").Value = "7.08"
$ws.Range("D27").Value = "0.999"
$ws.Range("D28").Value = "418.80"
$ws.Range("D31").Value = "160.97"
$ws.Range("D32").Value = "18.93"
$ws.Range("D34").Value = "17.79"
$ws.Range("D40").Value = "2.00"
$ws.Range("D41").Value = "3.33"
$ws.Range("D42").Value = "128.73"
$ws.Range("D45").Value = "0.554"
$ws.Range("D46").Value = "0.0911"
$ws.Range("D49").Value = "16.56"
$ws.Range("D51").Value = "0.0428"

# Reset the style back to Normal so no stray style index remains on the cell
foreach ($addr in $forceTextAddrs) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining updates (Volume % column, and Price column values that are not
# parsed as plain numbers, e.g. "67.702.26") can be set directly.
$ws.Range("D2").Value = "67.702.26"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.427.22"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("E9").Value = "  +8.62%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "67.610.03"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +6.13%  "
$ws.Range("E51").Value = "  +1.95%  "
